$p = $ppt.ActivePresentation

# The duplicate "Conclusion of Mad Mints" slide (position 7) is being
# removed from the deck.
$p.Slides.Item(7).Delete()

# Every slide that used to come after it shifts up by one position, so
# its cached "<N>" slide-number field text needs to be decremented by one
# to stay in sync (firstSlideNum="0", so a slide at 1-based position P
# displays "P-1").
for ($i = 7; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.Type -eq 14 -and $shape.PlaceholderFormat.Type -eq 13) {
            $shape.TextFrame.TextRange.Text = [string]($i - 1)
        }
    }
}
